$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Unique special functions, such as the Wheel of Fortune", $true, $false, $false, $false, $false,
    $true, 1, $false, "Unique special features", 2)

$d.Content.Find.Execute(
    "Eye-catching graphics with nostalgic design", $true, $false, $false, $false, $false,
    $true, 1, $false, "Eye-catching graphics", 2)

$d.Content.Find.Execute(
    "Immersive gaming experience with perfect sound design", $true, $false, $false, $false, $false,
    $true, 1, $false, "Immersive sound effects", 2)

$d.Content.Find.Execute(
    "Offers a range of in-game prizes and jackpots", $true, $false, $false, $false, $false,
    $true, 1, $false, "Range of prizes and jackpots", 2)

$d.Content.Find.Execute(
    "Graphics are simple and might not be appealing to everyone", $true, $false, $false, $false, $false,
    $true, 1, $false, "Simple graphics", 2)

$d.Content.Find.Execute(
    "Might feel repetitive after a while for some players", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited to fans of the television program", 2)

$d.Content.Find.Execute(
    "Read our review of MegaJackpots Wheel of Fortune On Air and play it for free. Enjoy unique special functions, immersive sound, and big jackpots.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of MegaJackpots Wheel of Fortune On Air and play for free.", 2)
